$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Introduction paragraph
Replace-Text "less battery, runtime memory" "less battery, and runtime memory"

# 2. History paragraph
Replace-Text "being introduced as beta in Android KitKat" "being introduced as a beta in Android KitKat"
Replace-Text "finally became the integral part" "finally became an integral part"

# 3. Dalvik's Android Era paragraph
Replace-Text "Dalvik was developed by Dan Bornstein which has origins from a town in Iceland" "Dalvik was developed by Dan Bornstein and has origins in a town in Iceland"

# 4. Dalvik's Architecture paragraph
Replace-Text "at that moment of time and the rest" "at that moment and the rest"
Replace-Text "Android Mobile phones are limited whether it be battery life" "Android Mobile phones is limited whether it be the battery life"

# 5. ART's Android Era paragraph
Replace-Text "ART has totally replaced the DVM" "ART has replaced the DVM"
Replace-Text "improve the performance for High-End Applications" "improve the performance of High-End Applications"
Replace-Text "overcome the Drawbacks with Dalvik." "overcome the Drawbacks of Dalvik."

# 6. ART's Architecture paragraph
Replace-Text "compiles the DEX code on the time of installation" "compiles the DEX code at the time of installation"

# 7. Caption - Architecture of DVM and ART
Replace-Text "Architecture of DVM and ART" "The architecture of DVM and ART"

# 8. Benefits of ART paragraph
Replace-Text "With Dalvik the problem was" "With Dalvik, the problem was"
Replace-Text "a lot of CPU, battery and a lot of garbage data gets created" "a lot of CPU, and battery and a lot of garbage data get created"
Replace-Text "which leads to decrease in the performance" "which leads to a decrease in the performance"
Replace-Text "using ART over the Dalvik is that It has better" "using ART over the Dalvik are that It has better"
Replace-Text "redundancy and it is AOT compiler" "redundancy and it is an AOT compiler"
Replace-Text "hence less Application Start up time" "hence less Application start up time"

# 9. Drawbacks of ART paragraph
Replace-Text "native code must be stored hence it requires larger physical memory" "native code must be stored it requires larger physical memory"

# 10. AOT? or JIT? paragraph
Replace-Text "So Why to use only AOT or JIT why not to use Both AOT and JIT" "So Why use only AOT or JIT why not use Both AOT and JIT"
Replace-Text "other’s drawbacks. Use of ART JIT compiler saves" "other’s drawbacks. The use of the ART JIT compiler saves"

# 11. Caption - Architecture of ART JIT
Replace-Text "Architecture of ART JIT " "The architecture of ART JIT "
